# Update the lattice-multiplication exercise table: every cell's five
# lines (problem header, the two multiplier digits, the "----" rule,
# and the two partial-product stub lines) are replaced in place, cell
# by cell, row by row. Table shape (5 rows x 3 cols) is unchanged.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11   # vertical-tab == <w:br/> line break within a run

# Each entry: row, col, then the 5 text lines for that cell.
$cells = @(
    @(1,1,"35 x 86","  8    6","  ----","3|    |","5|    |"),
    @(1,2,"16 x 62","  6    2","  ----","1|    |","6|    |"),
    @(1,3,"54 x 12","  1    2","  ----","5|    |","4|    |"),

    @(2,1,"83 x 83","  8    3","  ----","8|    |","3|    |"),
    @(2,2,"50 x 67","  6    7","  ----","5|    |","0|    |"),
    @(2,3,"52 x 60","  6    0","  ----","5|    |","2|    |"),

    @(3,1,"79 x 90","  9    0","  ----","7|    |","9|    |"),
    @(3,2,"22 x 71","  7    1","  ----","2|    |","2|    |"),
    @(3,3,"46 x 29","  2    9","  ----","4|    |","6|    |"),

    @(4,1,"13 x 85","  8    5","  ----","1|    |","3|    |"),
    @(4,2,"73 x 24","  2    4","  ----","7|    |","3|    |"),
    @(4,3,"15 x 90","  9    0","  ----","1|    |","5|    |"),

    @(5,1,"83 x 67","  6    7","  ----","8|    |","3|    |"),
    @(5,2,"11 x 22","  2    2","  ----","1|    |","1|    |"),
    @(5,3,"69 x 10","  1    0","  ----","6|    |","9|    |")
)

foreach ($entry in $cells) {
    $row = $entry[0]
    $col = $entry[1]
    $lines = $entry[2..6]
    $newText = [string]::Join($vt, $lines)
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
